$data = New-Object 'object[,]' 12,20
$data[0,0] = "ECs"
$data[0,1] = "Ucn2"
$data[0,2] = "Crhr2"
$data[0,3] = "ECs"
$data[0,4] = 2.0
$data[0,5] = 1.0
$data[0,6] = 0.772903
$data[0,7] = 1.545806
$data[0,8] = 0.05922444023038517
$data[0,9] = 0.05686984552343932
$data[0,10] = 2.0
$data[0,11] = 1.0
$data[0,12] = 1.2696695
$data[0,13] = 2.539339
$data[0,14] = 0.2673143181973693
$data[0,15] = 0.2673143181973693
$data[0,16] = 0.9813313655585
$data[0,17] = 3.925325462234
$data[0,18] = 0.01583154086080626
$data[0,19] = 0.0152021239820879
$data[1,0] = "ECs"
$data[1,1] = "Ucn2"
$data[1,2] = "Crhr2"
$data[1,3] = "MuSCs"
$data[1,4] = 2.0
$data[1,5] = 1.0
$data[1,6] = 0.772903
$data[1,7] = 1.545806
$data[1,8] = 0.05922444023038517
$data[1,9] = 0.05686984552343932
$data[1,10] = 2.0
$data[1,11] = 1.0
$data[1,12] = 3.4800555
$data[1,13] = 6.960110999999999
$data[1,14] = 0.7326856818026307
$data[1,15] = 0.7326856818026307
$data[1,16] = 2.6897453361165
$data[1,17] = 10.758981344466
$data[1,18] = 0.04339289936957891
$data[1,19] = 0.04166772154135143
$data[2,0] = "FAPs"
$data[2,1] = "Ucn2"
$data[2,2] = "Crhr2"
$data[2,3] = "ECs"
$data[2,4] = 3.0
$data[2,5] = 1.0
$data[2,6] = 0.2182836666666667
$data[2,7] = 0.6548510000000001
$data[2,8] = 0.01672619716803962
$data[2,9] = 0.02409181696207012
$data[2,10] = 2.0
$data[2,11] = 1.0
$data[2,12] = 1.2696695
$data[2,13] = 2.539339
$data[2,14] = 0.2673143181973693
$data[2,15] = 0.2673143181973693
$data[2,16] = 0.2771481139148333
$data[2,17] = 1.662888683489
$data[2,18] = 0.00447115199200928
$data[2,19] = 0.006440087625351592
$data[3,0] = "FAPs"
$data[3,1] = "Ucn2"
$data[3,2] = "Crhr2"
$data[3,3] = "MuSCs"
$data[3,4] = 3.0
$data[3,5] = 1.0
$data[3,6] = 0.2182836666666667
$data[3,7] = 0.6548510000000001
$data[3,8] = 0.01672619716803962
$data[3,9] = 0.02409181696207012
$data[3,10] = 2.0
$data[3,11] = 1.0
$data[3,12] = 3.4800555
$data[3,13] = 6.960110999999999
$data[3,14] = 0.7326856818026307
$data[3,15] = 0.7326856818026307
$data[3,16] = 0.7596392747435
$data[3,17] = 4.557835648461
$data[3,18] = 0.01225504517603034
$data[3,19] = 0.01765172933671853
$data[4,0] = "Inflammatory-Mac"
$data[4,1] = "Ucn2"
$data[4,2] = "Crhr2"
$data[4,3] = "ECs"
$data[4,4] = 1.0
$data[4,5] = 0.3333333333333333
$data[4,6] = 0.1300953333333333
$data[4,7] = 0.390286
$data[4,8] = 0.009968680795975739
$data[4,9] = 0.01435853174975452
$data[4,10] = 2.0
$data[4,11] = 1.0
$data[4,12] = 1.2696695
$data[4,13] = 2.539339
$data[4,14] = 0.2673143181973693
$data[4,15] = 0.2673143181973693
$data[4,16] = 0.1651780768256667
$data[4,17] = 0.991068460954
$data[4,18] = 0.002664771110303464
$data[4,19] = 0.003838241125000911
$data[5,0] = "Inflammatory-Mac"
$data[5,1] = "Ucn2"
$data[5,2] = "Crhr2"
$data[5,3] = "MuSCs"
$data[5,4] = 1.0
$data[5,5] = 0.3333333333333333
$data[5,6] = 0.1300953333333333
$data[5,7] = 0.390286
$data[5,8] = 0.009968680795975739
$data[5,9] = 0.01435853174975452
$data[5,10] = 2.0
$data[5,11] = 1.0
$data[5,12] = 3.4800555
$data[5,13] = 6.960110999999999
$data[5,14] = 0.7326856818026307
$data[5,15] = 0.7326856818026307
$data[5,16] = 0.452738980291
$data[5,17] = 2.716433881746
$data[5,18] = 0.007303909685672276
$data[5,19] = 0.01052029062475361
$data[6,0] = "MuSCs"
$data[6,1] = "Ucn2"
$data[6,2] = "Crhr2"
$data[6,3] = "ECs"
$data[6,4] = 2.0
$data[6,5] = 1.0
$data[6,6] = 11.1968455
$data[6,7] = 22.393691
$data[6,8] = 0.857969120424694
$data[6,9] = 0.8238587169862411
$data[6,10] = 2.0
$data[6,11] = 1.0
$data[6,12] = 1.2696695
$data[6,13] = 2.539339
$data[6,14] = 0.2673143181973693
$data[6,15] = 0.2673143181973693
$data[6,16] = 14.21629322756225
$data[6,17] = 56.865172910249
$data[6,18] = 0.2293474304607237
$data[6,19] = 0.2202292312221365
$data[7,0] = "MuSCs"
$data[7,1] = "Ucn2"
$data[7,2] = "Crhr2"
$data[7,3] = "MuSCs"
$data[7,4] = 2.0
$data[7,5] = 1.0
$data[7,6] = 11.1968455
$data[7,7] = 22.393691
$data[7,8] = 0.857969120424694
$data[7,9] = 0.8238587169862411
$data[7,10] = 2.0
$data[7,11] = 1.0
$data[7,12] = 3.4800555
$data[7,13] = 6.960110999999999
$data[7,14] = 0.7326856818026307
$data[7,15] = 0.7326856818026307
$data[7,16] = 38.96564376492525
$data[7,17] = 155.862575059701
$data[7,18] = 0.6286216899639703
$data[7,19] = 0.6036294857641047
$data[8,0] = "Neutrophils"
$data[8,1] = "Ucn2"
$data[8,2] = "Crhr2"
$data[8,3] = "ECs"
$data[8,4] = 2.0
$data[8,5] = 0.6666666666666666
$data[8,6] = 0.4889146666666667
$data[8,7] = 1.466744
$data[8,8] = 0.03746355940364922
$data[8,9] = 0.05396117281368521
$data[8,10] = 2.0
$data[8,11] = 1.0
$data[8,12] = 1.2696695
$data[8,13] = 2.539339
$data[8,14] = 0.2673143181973693
$data[8,15] = 0.2673143181973693
$data[8,16] = 0.6207600403693333
$data[8,17] = 3.724560242216
$data[8,18] = 0.01001454583923314
$data[8,19] = 0.01442459411982068
$data[9,0] = "Neutrophils"
$data[9,1] = "Ucn2"
$data[9,2] = "Crhr2"
$data[9,3] = "MuSCs"
$data[9,4] = 2.0
$data[9,5] = 0.6666666666666666
$data[9,6] = 0.4889146666666667
$data[9,7] = 1.466744
$data[9,8] = 0.03746355940364922
$data[9,9] = 0.05396117281368521
$data[9,10] = 2.0
$data[9,11] = 1.0
$data[9,12] = 3.4800555
$data[9,13] = 6.960110999999999
$data[9,14] = 0.7326856818026307
$data[9,15] = 0.7326856818026307
$data[9,16] = 1.701450174764
$data[9,17] = 10.208701048584
$data[9,18] = 0.02744901356441608
$data[9,19] = 0.03953657869386453
$data[10,0] = "Resolving-Mac"
$data[10,1] = "Ucn2"
$data[10,2] = "Crhr2"
$data[10,3] = "ECs"
$data[10,4] = 3.0
$data[10,5] = 1.0
$data[10,6] = 0.243364
$data[10,7] = 0.730092
$data[10,8] = 0.01864800197725647
$data[10,9] = 0.02685991596480985
$data[10,10] = 2.0
$data[10,11] = 1.0
$data[10,12] = 1.2696695
$data[10,13] = 2.539339
$data[10,14] = 0.2673143181973693
$data[10,15] = 0.2673143181973693
$data[10,16] = 0.308991848198
$data[10,17] = 1.853951089188
$data[10,18] = 0.004984877934293509
$data[10,19] = 0.007180040122971782
$data[11,0] = "Resolving-Mac"
$data[11,1] = "Ucn2"
$data[11,2] = "Crhr2"
$data[11,3] = "MuSCs"
$data[11,4] = 3.0
$data[11,5] = 1.0
$data[11,6] = 0.243364
$data[11,7] = 0.730092
$data[11,8] = 0.01864800197725647
$data[11,9] = 0.02685991596480985
$data[11,10] = 2.0
$data[11,11] = 1.0
$data[11,12] = 3.4800555
$data[11,13] = 6.960110999999999
$data[11,14] = 0.7326856818026307
$data[11,15] = 0.7326856818026307
$data[11,16] = 0.846920226702
$data[11,17] = 5.081521360211999
$data[11,18] = 0.01366312404296296
$data[11,19] = 0.01967987584183807

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the extra rows (14-16) first so row 13 stays where it is before we overwrite it
$ws.Rows("14:16").Delete()

# Write the full updated data block A2:T13
$ws.Range("A2:T13").Value = $data
